$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the grade values for rows 15-20 (columns C and D), G column formulas will recompute automatically
$ws.Range("C15").Value = 9

$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 10

$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 10

$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 10

$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 10

$ws.Range("C20").Value = 10

# Update the view: scroll so row 7 is the top-left visible row, and select E11
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("E11").Select()
